# Appointments Module Update:
# Adds new test-data rows for the "Cancel Appointment" and
# "Verify video invites" scenarios to the MMH testdata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: APPOINTMENT_DETAILS_AFTER_CANCELLED
# (set column B/C first so the new shared strings are interned in the
# same order as the authored workbook: "Cancelled;..." before the KEY)
$ws.Range("B27").Value = "Cancelled;Dr Sam Entwistle;"
$ws.Range("C27").Value = "Cancelled;Dr Sam Entwistle;"
$ws.Range("A27").Value = "APPOINTMENT_DETAILS_AFTER_CANCELLED"

# Row 28 intentionally left blank (matches the blank-separator pattern
# used throughout the sheet).

# Row 29: VIDEO_ICONS
$ws.Range("A29").Value = "VIDEO_ICONS"
$ws.Range("B29").Value = "Minimize;Close;Settings;Mute;Share Now;Video Off"
$ws.Range("C29").Value = "Minimize;Close;Settings;Mute;Share Now;Video Off"

# Match the resulting selection/active cell left behind in the sheet.
$ws.Range("A29").Select() | Out-Null
